$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.404.36'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.693.60'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.24'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5493'
$ws.Range('E6').Value = '  +4.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2735'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06472'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.05'
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07675'
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('D12').Value = '1.695.29'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.556'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5856'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008385'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.49'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').Value = '26.456.37'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.010'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.65'
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.257'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.31'
$ws.Range('E24').Value = '  +3.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1332'
$ws.Range('E25').Value = '  +7.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.918'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.80'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06303'
$ws.Range('E28').Value = '  -5.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.397'
$ws.Range('E29').Value = '  +3.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.331'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.606'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.611'
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.686'
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6151'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.196'
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = '1.119.83'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8865'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.96'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.54'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.232'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.008'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.108'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4305'
$ws.Range('E51').Value = '  +0.11%  '
